# "Generate Report for Handoff" - mark e2e\b.md as ready for a new handoff
# round because its last handback is stale (a newer source revision exists).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/7a5e2a4c7c3fd3c920a97d178d27d467857f751a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/13f2d288a3a8dc3b109257cc01867f539a295d7f/e2e/b.md."

# --- Overview sheet: row 3 is the b.md file ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-13 22:48:23"

# --- zh-cn sheet: row 3 is the b.md file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-13 22:48:15"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the b.md file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-13 22:48:23"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
